$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 120182.6963248574
$ws.Range("E4").Value = 0.001013526776580606
$ws.Range("F4").Value = 0.2143536725542664
$ws.Range("G4").Value = -0.9571585894267705
$ws.Range("H4").Value = 10.68364220396225

$ws.Range("D5").Value = 120983.0067852185
$ws.Range("E5").Value = -0.01088978681485037
$ws.Range("F5").Value = 0.2445175544675988
$ws.Range("G5").Value = -1.370914664415405
$ws.Range("H5").Value = 13.7504751161537

$ws.Range("D7").Value = 122293.3437093496
$ws.Range("E7").Value = -0.02376661751950196
$ws.Range("F7").Value = 0.260494992602112
$ws.Range("G7").Value = -1.129483318729288
$ws.Range("H7").Value = 10.71342127783571

$ws.Range("D8").Value = 123937.0595087343
$ws.Range("E8").Value = -0.0390933828783337
$ws.Range("F8").Value = 0.2129417120102251
$ws.Range("G8").Value = -0.7623385489216303
$ws.Range("H8").Value = 6.659797727762894

$ws.Range("D9").Value = 125253.2435387082
$ws.Range("E9").Value = -0.07266281220699891
$ws.Range("F9").Value = 0.3310479510830294
$ws.Range("G9").Value = -1.542482876485479
$ws.Range("H9").Value = 9.932864118809061

$ws.Range("D10").Value = 126820.8199628138
$ws.Range("E10").Value = -0.1110221762133587
$ws.Range("F10").Value = 0.4429786770403688
$ws.Range("G10").Value = -1.947676407540862
$ws.Range("H10").Value = 9.986871299635201

$ws.Range("D11").Value = 128939.2712006311
$ws.Range("E11").Value = -0.1899320620052032
$ws.Range("F11").Value = 0.7873168261214663
$ws.Range("G11").Value = -2.661611395256252
$ws.Range("H11").Value = 13.12153382279394

$ws.Range("D17").Value = 119391.4152366399
$ws.Range("E17").Value = 0.06438063257754784
$ws.Range("F17").Value = 0.1611339399763616
$ws.Range("G17").Value = -1.722164190215021
$ws.Range("H17").Value = 13.61415250937092

$ws.Range("D18").Value = 119539.1105384298
$ws.Range("E18").Value = 0.04744025577969273
$ws.Range("F18").Value = 0.1565464439957156
$ws.Range("G18").Value = -0.6389219314865378
$ws.Range("H18").Value = 9.62109017717134

$ws.Range("D19").Value = 119576.4072627097
$ws.Range("E19").Value = 0.03448563976614866
$ws.Range("F19").Value = 0.1591395647225821
$ws.Range("G19").Value = -0.3104426309802923
$ws.Range("H19").Value = 6.674249834399131

$ws.Range("D20").Value = 119616.1874060205
$ws.Range("E20").Value = 0.02367896888974719
$ws.Range("F20").Value = 0.1711197132699751
$ws.Range("G20").Value = -0.2379160475045328
$ws.Range("H20").Value = 5.732725964501296

